$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.092.37"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.830.13"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.58%  "
$ws.Range("D5").Value = "'312.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.71%  "
$ws.Range("E6").Value = "  +0.62%  "
$ws.Range("D7").Value = "'0.4705"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.3682"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").Value = "'0.07380"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "'0.8796"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.841.98"
$ws.Range("E12").Value = "  -1.12%  "
$ws.Range("D13").Value = "'0.07292"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.23%  "
$ws.Range("D14").Value = "'5.459"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.81%  "
$ws.Range("D15").Value = "'92.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("D16").Value = "'6.544"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "'1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("D18").Value = "'0.000008765"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "
$ws.Range("D19").Value = "'1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "27.115.16"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'5.308"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.77%  "
$ws.Range("E23").Value = "  +1.07%  "
$ws.Range("D24").Value = "2.052.22"
$ws.Range("E24").Value = "  -2.24%  "
$ws.Range("D25").Value = "'1.895"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'152.02"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").Value = "'18.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.10%  "
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "'5.266"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.47%  "
$ws.Range("D30").Value = "'117.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.41%  "
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "'0.7583"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").Value = "'1.166"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").Value = "'4.531"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "'2.928"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("E36").Value = "  +0.87%  "
$ws.Range("D37").Value = "'1.102"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").Value = "'0.05327"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.28%  "
$ws.Range("D39").Value = "'0.01956"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'2.991"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").Value = "'2.418"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.02%  "
$ws.Range("D42").Value = "'7.278"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").Value = "'0.5334"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "
$ws.Range("D44").Value = "'0.1661"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").Value = "'8.540"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.78%  "
$ws.Range("D46").Value = "'0.4939"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "'10.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.45%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "'1.668"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.62%  "
$ws.Range("D50").Value = "'103.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'0.06308"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.17%  "
